$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("M21").Value = "icon_reptile_0"
$ws.Range("M22").Value = "coin_visual"
$ws.Range("M23").Value = "gem_visual"
$ws.Range("N19").Value = "TID_XPROMO_REWARD_EGG_BETTER"
$ws.Range("N20").Value = "TID_XPROMO_REWARD_PET_33"
$ws.Range("N21").Value = " TID_XPROMO_REWARD_DRAGON_REPTILE"
$ws.Range("N22").Value = "FE_POPUP_IAP_COINS_LOWERCASE"
$ws.Range("N23").Value = "FE_POPUP_IAP_GEMS_LOWERCASE"
$ws.Range("N24").Value = "TID_XPROMO_REWARD_PET_33"

$ws.Columns.Item(14).ColumnWidth = 37.6

$ws.Range("N27").Select()
